# Applies the two logical changes described by the commit:
#   1. Bump the cached "datetimeFigureOut" footer-date field from
#      17/01/2018 to 22/01/2018 everywhere it is cached (the slide
#      master and every slide layout).
#   2. Mark two more to-do bullets as done (single strikethrough) in the
#      "CaixaDeTexto 9" shape on slide 1 - the "- desenhar tela para
#      inicio de jogo" line and the "- add save/load logic to menus"
#      line - without touching the leading tab run that starts each
#      paragraph.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: Presentation.SlideMaster + every
#    CustomLayout hanging off it each own their own cached copy of the
#    "dt" (date) placeholder.
# ---------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = "22/01/2018"
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateShapes $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Strikethrough the two newly-finished bullets on slide 1.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$shape = $slide1.Shapes.Item(7)   # "CaixaDeTexto 9"
$tr = $shape.TextFrame.TextRange

# "	- desenhar tela para inicio de jogo" -> paragraph 6.
# Skip the leading tab run; only strike "- desenhar tela para inicio de jogo".
$line1 = $tr.Paragraphs(6, 1)
$tabLen = 1
$rest = $tr.Characters($line1.Start + $tabLen, $line1.Length - $tabLen - 1)
$rest.Font.Strikethrough = -1

# "	- add save/load logic to menus" -> paragraph 8.
# Skip the leading tab run; only strike "- add save/load logic to menus".
$line2 = $tr.Paragraphs(8, 1)
$rest2 = $tr.Characters($line2.Start + $tabLen, $line2.Length - $tabLen - 1)
$rest2.Font.Strikethrough = -1
